$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '25.833.31'
$ws.Range("E2").Value = '  -0.33%  '
$ws.Range("D3").Value = '1.636.42'
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.003'
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = '  +0.24%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '215.27'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -0.64%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.5044'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -0.55%  '
$ws.Range("E7").Value = '  +0.25%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.2572'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  -0.26%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.06415'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +0.90%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '20.18'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +2.78%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.07791'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +0.58%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '4.290'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +0.72%  '
$ws.Range("D13").Value = '1.861.83'
$ws.Range("E13").Value = '  +0.00%  '
$ws.Range("D14").Value = '1.635.54'
$ws.Range("E14").Value = '  +0.17%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.5605'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +1.63%  '
$ws.Range("D16").Value = '0.0₅7636'
$ws.Range("E16").Value = '  -0.95%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '63.05'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  -1.44%  '
$ws.Range("D18").Value = '25.857.24'
$ws.Range("E18").Value = '  -0.30%  '
$ws.Range("E19").Value = '  +0.19%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '194.62'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -0.05%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '4.331'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -2.63%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '9.885'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -0.19%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '6.099'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +0.64%  '
$ws.Range("E24").Value = '  +0.20%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '1.780'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -6.88%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '140.44'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -1.30%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '0.1254'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +1.39%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '6.810'
$ws.Range("D28").Style = "Normal"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '15.48'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -0.55%  '
$ws.Range("E30").Value = '  -0.45%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '0.04904'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +0.48%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '3.302'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +1.50%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '3.235'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +1.46%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '1.571'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +1.81%  '
$ws.Range("E35").Value = '  +0.11%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.9033'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  -0.33%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '2.577'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +0.37%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.5534'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +0.68%  '
$ws.Range("D39").Value = '1.125.88'
$ws.Range("E39").Value = '  +0.30%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.01558'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +0.09%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '1.002'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +0.09%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '5.517'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -1.01%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.7992'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -1.08%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '98.11'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +0.59%  '
$ws.Range("D45").Value = '1.772.10'
$ws.Range("E45").Value = '  -0.25%  '
$ws.Range("E46").Value = '  -7.00%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '55.46'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +1.03%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.4266'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -4.26%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '7.724'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +2.70%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.05035'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -2.18%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '1.001'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +0.54%  '
